# Auto-generated Excel COM-interop script applying the cryptos.xlsx data refresh
# (diff of Fri Mar  8 15:07:46 UTC 2024 GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.479.20"
$ws.Range("E2").Value = "  +1.94%  "
$ws.Range("D3").Value = "3.963.45"
$ws.Range("E3").Value = "  +3.86%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("D5").Value = "'483.50"
$ws.Range("E5").Value = "  +8.15%  "
$ws.Range("D6").Value = "'151.21"
$ws.Range("E6").Value = "  +2.88%  "
$ws.Range("D7").Value = "'0.626"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "'0.733"
$ws.Range("E9").Value = "  -1.03%  "
$ws.Range("D10").Value = "'0.169"
$ws.Range("E10").Value = "  +7.91%  "
$ws.Range("D11").Value = "'0.0000355"
$ws.Range("E11").Value = "  +9.94%  "
$ws.Range("D12").Value = "'43.61"
$ws.Range("E12").Value = "  -1.08%  "
$ws.Range("D13").Value = "4.605.60"
$ws.Range("E13").Value = "  +4.15%  "
$ws.Range("D14").Value = "'10.47"
$ws.Range("E14").Value = "  +1.32%  "
$ws.Range("D15").Value = "'14.79"
$ws.Range("E15").Value = "  -0.30%  "
$ws.Range("D16").Value = "3.955.71"
$ws.Range("E16").Value = "  +4.67%  "
$ws.Range("E17").Value = "  +0.08%  "
$ws.Range("D18").Value = "'20.02"
$ws.Range("E18").Value = "  +0.07%  "
$ws.Range("D19").Value = "'1.14"
$ws.Range("E19").Value = "  -0.74%  "
$ws.Range("D20").Value = "68.590.17"
$ws.Range("E20").Value = "  +2.04%  "
$ws.Range("D21").Value = "'437.30"
$ws.Range("E21").Value = "  +3.55%  "
$ws.Range("D22").Value = "'3.36"
$ws.Range("E22").Value = "  +2.96%  "
$ws.Range("D23").Value = "'14.39"
$ws.Range("E23").Value = "  -2.02%  "
$ws.Range("D24").Value = "'88.10"
$ws.Range("E24").Value = "  +1.60%  "
$ws.Range("D25").Value = "'3.62"
$ws.Range("E25").Value = "  +6.19%  "
$ws.Range("D26").Value = "'38.48"
$ws.Range("E26").Value = "  +3.05%  "
$ws.Range("D27").Value = "'10.08"
$ws.Range("E27").Value = "  +2.75%  "
$ws.Range("E28").Value = "  +3.80%  "
$ws.Range("D29").Value = "'727.72"
$ws.Range("E29").Value = "  -1.09%  "
$ws.Range("D30").Value = "'13.24"
$ws.Range("E30").Value = "  -3.31%  "
$ws.Range("D31").Value = "'0.127"
$ws.Range("E31").Value = "  -4.96%  "
$ws.Range("D32").Value = "'2.83"
$ws.Range("E32").Value = "  +3.36%  "
$ws.Range("B33").Value = "InjectiveProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D33").Value = "'42.67"
$ws.Range("E33").Value = "  -1.53%  "
$ws.Range("B34").Value = "PEPE"
$ws.Range("C34").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D34").Value = "0.0₃0859"
$ws.Range("E34").Value = "  +26.50%  "
$ws.Range("D35").Value = "'59.97"
$ws.Range("E35").Value = "  +5.64%  "
$ws.Range("E36").Value = "  -4.83%  "
$ws.Range("D37").Value = "'0.998"
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("D38").Value = "'5.38"
$ws.Range("E38").Value = "  -2.68%  "
$ws.Range("D39").Value = "'0.0474"
$ws.Range("E39").Value = "  -0.84%  "
$ws.Range("D40").Value = "'3.04"
$ws.Range("E40").Value = "  +4.93%  "
$ws.Range("D41").Value = "'2.89"
$ws.Range("E41").Value = "  +8.28%  "
$ws.Range("B42").Value = "Fetch.AI"
$ws.Range("C42").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D42").Value = "'2.58"
$ws.Range("E42").Value = "  +2.08%  "
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").Value = "'0.141"
$ws.Range("E43").Value = "  +0.39%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  +0.12%  "
$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").Value = "'2.21"
$ws.Range("E45").Value = "  +4.58%  "
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").Value = "'0.335"
$ws.Range("E46").Value = "  -1.77%  "
$ws.Range("D47").Value = "'3.42"
$ws.Range("E47").Value = "  -0.15%  "
$ws.Range("D48").Value = "'3.24"
$ws.Range("E48").Value = "  -0.65%  "
$ws.Range("D49").Value = "'148.40"
$ws.Range("E49").Value = "  +1.21%  "
$ws.Range("E50").Value = "  +0.27%  "
$ws.Range("D51").Value = "'25.03"
$ws.Range("E51").Value = "  -1.03%  "
